# Rename the embedded logo pictures that live in the document's
# headers/footers:
#
#   - Pearson logo (appears in the footers) : image1.png -> image2.png
#   - BTec logo    (appears in the headers) : image2.jpg -> image1.jpg
#
# Word exposes a picture's OOXML name (wp:docPr/@name) through the
# InlineShape.Name property, so the rename is done through the object
# model rather than by touching raw XML.

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    # --- Headers: BTec_Logo-Orange --------------------------------------
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    # Re-fetch right before writing - some header/footer
                    # shape handles go stale otherwise.
                    $fresh = $shp.Range.InlineShapes.Item(1)
                    $fresh.Name = "image1.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson logo -------------------------------------------
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $fresh = $shp.Range.InlineShapes.Item(1)
                    $fresh.Name = "image2.png"
                }
            }
        }
    }
}
